# Update the LR-pairs sheet (Ly9-Ly9) with the re-run TPM results.
# A new "ECs" cluster was added to the sending/target cluster set, so the
# cross-product of clusters now yields 9 data rows (3 x 3) instead of 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('ECs', 'Ly9', 'Ly9', 'ECs', '1', '0.3333333333333333', '0.0435', '0.1305', '0.001251028107840694', '0.001251028107840694', '1', '0.3333333333333333', '0.0435', '0.1305', '0.001251028107840694', '0.001251028107840694', '0.00189225', '0.01703025', '1.565071326607467E-06', '1.565071326607466E-06'),
    @('ECs', 'Ly9', 'Ly9', 'MuSCs', '1', '0.3333333333333333', '0.0435', '0.1305', '0.001251028107840694', '0.001251028107840694', '2', '0.6666666666666666', '0.02618033333333333', '0.078541', '0.0007529271924744515', '0.0007529271924744514', '0.0011388445', '0.0102496005', '9.41933080943119E-07', '9.419330809431186E-07'),
    @('ECs', 'Ly9', 'Ly9', 'Resolving-Mac', '1', '0.3333333333333333', '0.0435', '0.1305', '0.001251028107840694', '0.001251028107840694', '3', '1', '34.70172066666667', '104.105162', '0.9979960446996849', '0.9979960446996848', '1.509524849', '13.585723641', '0.001248521103433143', '0.001248521103433143'),
    @('MuSCs', 'Ly9', 'Ly9', 'ECs', '2', '0.6666666666666666', '0.02618033333333333', '0.078541', '0.0007529271924744515', '0.0007529271924744514', '1', '0.3333333333333333', '0.0435', '0.1305', '0.001251028107840694', '0.001251028107840694', '0.0011388445', '0.0102496005', '9.41933080943119E-07', '9.419330809431186E-07'),
    @('MuSCs', 'Ly9', 'Ly9', 'MuSCs', '2', '0.6666666666666666', '0.02618033333333333', '0.078541', '0.0007529271924744515', '0.0007529271924744514', '2', '0.6666666666666666', '0.02618033333333333', '0.078541', '0.0007529271924744515', '0.0007529271924744514', '0.0006854098534444445', '0.006168688681', '5.668993571674598E-07', '5.668993571674596E-07'),
    @('MuSCs', 'Ly9', 'Ly9', 'Resolving-Mac', '2', '0.6666666666666666', '0.02618033333333333', '0.078541', '0.0007529271924744515', '0.0007529271924744514', '3', '1', '34.70172066666667', '104.105162', '0.9979960446996849', '0.9979960446996848', '0.9085026142935555', '8.176523528642001', '0.000751418360036341', '0.0007514183600363408'),
    @('Resolving-Mac', 'Ly9', 'Ly9', 'ECs', '3', '1', '34.70172066666667', '104.105162', '0.9979960446996849', '0.9979960446996848', '1', '0.3333333333333333', '0.0435', '0.1305', '0.001251028107840694', '0.001251028107840694', '1.509524849', '13.585723641', '0.001248521103433143', '0.001248521103433143'),
    @('Resolving-Mac', 'Ly9', 'Ly9', 'MuSCs', '3', '1', '34.70172066666667', '104.105162', '0.9979960446996849', '0.9979960446996848', '2', '0.6666666666666666', '0.02618033333333333', '0.078541', '0.0007529271924744515', '0.0007529271924744514', '0.9085026142935555', '8.176523528642001', '0.000751418360036341', '0.0007514183600363408'),
    @('Resolving-Mac', 'Ly9', 'Ly9', 'Resolving-Mac', '3', '1', '34.70172066666667', '104.105162', '0.9979960446996849', '0.9979960446996848', '3', '1', '34.70172066666667', '104.105162', '0.9979960446996849', '0.9979960446996848', '1204.209417227361', '10837.88475504625', '0.9959961052362155', '0.9959961052362153')

)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $row.Length; $j++) {
        $col = $j + 1
        if ($j -lt 4) {
            $ws.Cells.Item($r, $col).Value = $row[$j]
        } else {
            $ws.Cells.Item($r, $col).Value = [double]$row[$j]
        }
    }
}
